$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 8.841467
$ws.Cells.Item(2, 8).Value = 26.524401
$ws.Cells.Item(2, 9).Value = 0.5917001192060068
$ws.Cells.Item(2, 10).Value = 0.5917001192060067
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.01631833333333333
$ws.Cells.Item(2, 14).Value = 0.048955
$ws.Cells.Item(2, 15).Value = 0.001076315602073535
$ws.Cells.Item(2, 16).Value = 0.001076315602073535
$ws.Cells.Item(2, 17).Value = 0.1442780056616667
$ws.Cells.Item(2, 18).Value = 1.298502050955
$ws.Cells.Item(2, 19).Value = 0.0006368560700501958
$ws.Cells.Item(2, 20).Value = 0.0006368560700501956

$ws.Cells.Item(3, 7).Value = 8.841467
$ws.Cells.Item(3, 8).Value = 26.524401
$ws.Cells.Item(3, 9).Value = 0.5917001192060068
$ws.Cells.Item(3, 10).Value = 0.5917001192060067
$ws.Cells.Item(3, 15).Value = 0.7730166590262294
$ws.Cells.Item(3, 16).Value = 0.7730166590262293
$ws.Cells.Item(3, 17).Value = 103.62137433731
$ws.Cells.Item(3, 18).Value = 932.59236903579
$ws.Cells.Item(3, 19).Value = 0.4573940492940491
$ws.Cells.Item(3, 20).Value = 0.4573940492940489

$ws.Cells.Item(4, 7).Value = 8.841467
$ws.Cells.Item(4, 8).Value = 26.524401
$ws.Cells.Item(4, 9).Value = 0.5917001192060068
$ws.Cells.Item(4, 10).Value = 0.5917001192060067
$ws.Cells.Item(4, 13).Value = 3.425042
$ws.Cells.Item(4, 14).Value = 10.275126
$ws.Cells.Item(4, 15).Value = 0.2259070253716972
$ws.Cells.Item(4, 16).Value = 0.2259070253716972
$ws.Cells.Item(4, 17).Value = 30.282395816614
$ws.Cells.Item(4, 18).Value = 272.541562349526
$ws.Cells.Item(4, 19).Value = 0.1336692138419076
$ws.Cells.Item(4, 20).Value = 0.1336692138419076

$ws.Cells.Item(5, 9).Value = 0.2746155987184545
$ws.Cells.Item(5, 10).Value = 0.2746155987184545
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.01631833333333333
$ws.Cells.Item(5, 14).Value = 0.048955
$ws.Cells.Item(5, 15).Value = 0.001076315602073535
$ws.Cells.Item(5, 16).Value = 0.001076315602073535
$ws.Cells.Item(5, 17).Value = 0.06696126909666666
$ws.Cells.Item(5, 18).Value = 0.6026514218699999
$ws.Cells.Item(5, 19).Value = 0.0002955730534734377
$ws.Cells.Item(5, 20).Value = 0.0002955730534734376

$ws.Cells.Item(6, 9).Value = 0.2746155987184545
$ws.Cells.Item(6, 10).Value = 0.2746155987184545
$ws.Cells.Item(6, 15).Value = 0.7730166590262294
$ws.Cells.Item(6, 16).Value = 0.7730166590262293
$ws.Cells.Item(6, 19).Value = 0.2122824326378274
$ws.Cells.Item(6, 20).Value = 0.2122824326378273

$ws.Cells.Item(7, 9).Value = 0.2746155987184545
$ws.Cells.Item(7, 10).Value = 0.2746155987184545
$ws.Cells.Item(7, 13).Value = 3.425042
$ws.Cells.Item(7, 14).Value = 10.275126
$ws.Cells.Item(7, 15).Value = 0.2259070253716972
$ws.Cells.Item(7, 16).Value = 0.2259070253716972
$ws.Cells.Item(7, 17).Value = 14.054447494396
$ws.Cells.Item(7, 18).Value = 126.490027449564
$ws.Cells.Item(7, 19).Value = 0.06203759302715371
$ws.Cells.Item(7, 20).Value = 0.0620375930271537

$ws.Cells.Item(8, 7).Value = 1.997574666666667
$ws.Cells.Item(8, 8).Value = 5.992724
$ws.Cells.Item(8, 9).Value = 0.1336842820755386
$ws.Cells.Item(8, 10).Value = 0.1336842820755386
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.01631833333333333
$ws.Cells.Item(8, 14).Value = 0.048955
$ws.Cells.Item(8, 15).Value = 0.001076315602073535
$ws.Cells.Item(8, 16).Value = 0.001076315602073535
$ws.Cells.Item(8, 17).Value = 0.03259708926888889
$ws.Cells.Item(8, 18).Value = 0.29337380342
$ws.Cells.Item(8, 19).Value = 0.0001438864785499016
$ws.Cells.Item(8, 20).Value = 0.0001438864785499016

$ws.Cells.Item(9, 7).Value = 1.997574666666667
$ws.Cells.Item(9, 8).Value = 5.992724
$ws.Cells.Item(9, 9).Value = 0.1336842820755386
$ws.Cells.Item(9, 10).Value = 0.1336842820755386
$ws.Cells.Item(9, 15).Value = 0.7730166590262294
$ws.Cells.Item(9, 16).Value = 0.7730166590262293
$ws.Cells.Item(9, 17).Value = 23.41143526310667
$ws.Cells.Item(9, 18).Value = 210.70291736796
$ws.Cells.Item(9, 19).Value = 0.1033401770943529
$ws.Cells.Item(9, 20).Value = 0.1033401770943529

$ws.Cells.Item(10, 7).Value = 1.997574666666667
$ws.Cells.Item(10, 8).Value = 5.992724
$ws.Cells.Item(10, 9).Value = 0.1336842820755386
$ws.Cells.Item(10, 10).Value = 0.1336842820755386
$ws.Cells.Item(10, 13).Value = 3.425042
$ws.Cells.Item(10, 14).Value = 10.275126
$ws.Cells.Item(10, 15).Value = 0.2259070253716972
$ws.Cells.Item(10, 16).Value = 0.2259070253716972
$ws.Cells.Item(10, 17).Value = 6.841777131469334
$ws.Cells.Item(10, 18).Value = 61.575994183224
$ws.Cells.Item(10, 19).Value = 0.03020021850263582
$ws.Cells.Item(10, 20).Value = 0.03020021850263582
